$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2373.3125
$ws.Range("J12").Value = 889
$ws.Range("L12").Value = 889
$ws.Range("N12").Value = -1229
$ws.Range("H17").Value = 3324.4348
$ws.Range("J17").Value = 3324.4348
$ws.Range("L17").Value = 9973.304400000001
$ws.Range("N17").Value = -10309.3044
$ws.Range("H19").Value = 1356.9231
$ws.Range("I19").Value = 1217.8889
$ws.Range("J19").Value = 1669.75
$ws.Range("K19").Value = 1217.8889
$ws.Range("L19").Value = 1669.75
$ws.Range("M19").Value = -1042.8889
$ws.Range("N19").Value = -2019.75
$ws.Range("H33").Value = 399.5
$ws.Range("I33").Value = 496.66666
$ws.Range("J33").Value = 224.6
$ws.Range("K33").Value = 496.66666
$ws.Range("L33").Value = 224.6
$ws.Range("M33").Value = -267.66666
$ws.Range("N33").Value = -682.6
$ws.Range("H43").Value = 999.6667
$ws.Range("I43").Value = 899.6667
$ws.Range("J43").Value = 1299.6666
$ws.Range("K43").Value = 899.6667
$ws.Range("L43").Value = 1299.6666
$ws.Range("M43").Value = -830.6667
$ws.Range("N43").Value = -1437.6666
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30924
$ws.Range("H98").Value = 1293.5834
$ws.Range("I98").Value = 684
$ws.Range("K98").Value = 684
$ws.Range("M98").Value = 814
$ws.Range("H99").Value = 344.5
$ws.Range("I99").Value = 344.5
$ws.Range("K99").Value = 1033.5
$ws.Range("M99").Value = 464.5
$ws.Range("H107").Value = 1028.4546
$ws.Range("I107").Value = 886.625
$ws.Range("J107").Value = 1406.6666
$ws.Range("K107").Value = 886.625
$ws.Range("L107").Value = 1406.6666
$ws.Range("M107").Value = 1033.375
$ws.Range("N107").Value = -5246.6666
$ws.Range("H116").Value = 4576.8
$ws.Range("I116").Value = 4496
$ws.Range("K116").Value = 4496
$ws.Range("M116").Value = -1054
$ws.Range("H122").Value = 1293.5834
$ws.Range("I122").Value = 684
$ws.Range("K122").Value = 2052
$ws.Range("M122").Value = 398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 212.73334
$ws.Range("I5").Value = 241.5
$ws.Range("J5").Value = 155.2
$ws.Range("K5").Value = 241.5
$ws.Range("L5").Value = 155.2
$ws.Range("M5").Value = -129.5
$ws.Range("N5").Value = -379.2
$ws.Range("H32").Value = 2080.389
$ws.Range("I32").Value = 1905.7354
$ws.Range("K32").Value = 1905.7354
$ws.Range("M32").Value = -1618.7354
$ws.Range("H122").Value = 6140
$ws.Range("I122").Value = 4663.3335
$ws.Range("K122").Value = 13990.0005
$ws.Range("M122").Value = -11540.0005
$ws.Range("H130").Value = 20391.8
$ws.Range("J130").Value = 19999
$ws.Range("L130").Value = 19999
$ws.Range("N130").Value = -30039

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 212.73334
$ws.Range("I4").Value = 241.5
$ws.Range("J4").Value = 155.2
$ws.Range("K4").Value = 241.5
$ws.Range("L4").Value = 155.2
$ws.Range("M4").Value = -126.5
$ws.Range("N4").Value = -385.2
$ws.Range("H22").Value = 262.83334
$ws.Range("I22").Value = 456.33334
$ws.Range("K22").Value = 456.33334
$ws.Range("M22").Value = -283.33334
$ws.Range("H134").Value = 5387.8623
$ws.Range("I134").Value = 5692.4546
$ws.Range("K134").Value = 17077.3638
$ws.Range("M134").Value = -14542.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1180.1666
$ws.Range("I16").Value = 1180.1666
$ws.Range("K16").Value = 1180.1666
$ws.Range("M16").Value = -893.1666
$ws.Range("H22").Value = 10000549
$ws.Range("I22").Value = 499
$ws.Range("K22").Value = 499
$ws.Range("M22").Value = -149
$ws.Range("H31").Value = 2193.4443
$ws.Range("J31").Value = 1586.2
$ws.Range("L31").Value = 1586.2
$ws.Range("N31").Value = -2176.2
$ws.Range("H34").Value = 2193.4443
$ws.Range("J34").Value = 1586.2
$ws.Range("L34").Value = 1586.2
$ws.Range("N34").Value = -1990.2
$ws.Range("H45").Value = 14999
$ws.Range("J45").Value = 14999
$ws.Range("L45").Value = 14999
$ws.Range("N45").Value = -16185
$ws.Range("H113").Value = 1180.1666
$ws.Range("I113").Value = 1180.1666
$ws.Range("K113").Value = 1180.1666
$ws.Range("M113").Value = 989.8334
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 10472.75
$ws.Range("J55").Value = 11464.333
$ws.Range("L55").Value = 34392.999
$ws.Range("N55").Value = -34746.999
$ws.Range("H80").Value = 6071.6
$ws.Range("I80").Value = 6057.3335
$ws.Range("K80").Value = 18172.0005
$ws.Range("M80").Value = -17236.0005
$ws.Range("H83").Value = 6071.6
$ws.Range("I83").Value = 6057.3335
$ws.Range("K83").Value = 54516.0015
$ws.Range("M83").Value = -49836.0015
$ws.Range("H107").Value = 881.5
$ws.Range("I107").Value = 847.25
$ws.Range("K107").Value = 2541.75
$ws.Range("M107").Value = -621.75
$ws.Range("H116").Value = 166999
$ws.Range("I116").Value = 166999
$ws.Range("K116").Value = 500997
$ws.Range("M116").Value = -497555
$ws.Range("H137").Value = 5324.6665
$ws.Range("I137").Value = 2983.3333
$ws.Range("K137").Value = 8949.999899999999
$ws.Range("M137").Value = -3849.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8961.6
$ws.Range("I70").Value = 8976
$ws.Range("J70").Value = 8940
$ws.Range("K70").Value = 8976
$ws.Range("L70").Value = 8940
$ws.Range("M70").Value = -8706
$ws.Range("N70").Value = -9480
$ws.Range("H73").Value = 8961.6
$ws.Range("I73").Value = 8976
$ws.Range("J73").Value = 8940
$ws.Range("K73").Value = 8976
$ws.Range("L73").Value = 8940
$ws.Range("M73").Value = -8040
$ws.Range("N73").Value = -10812
$ws.Range("H107").Value = 2953.2727
$ws.Range("I107").Value = 1812.7142
$ws.Range("J107").Value = 4949.25
$ws.Range("K107").Value = 1812.7142
$ws.Range("L107").Value = 4949.25
$ws.Range("M107").Value = 107.2858000000001
$ws.Range("N107").Value = -8789.25
$ws.Range("H122").Value = 2594.32
$ws.Range("I122").Value = 2471.476
$ws.Range("J122").Value = 3239.25
$ws.Range("K122").Value = 7414.428
$ws.Range("L122").Value = 9717.75
$ws.Range("M122").Value = -4964.428
$ws.Range("N122").Value = -14617.75
$ws.Range("H126").Value = 2892.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2173
$ws.Range("I7").Value = 2201.8948
$ws.Range("K7").Value = 2201.8948
$ws.Range("M7").Value = -2089.8948
$ws.Range("H40").Value = 1797.8
$ws.Range("I40").Value = 1797.8
$ws.Range("K40").Value = 1797.8
$ws.Range("M40").Value = -1661.8
$ws.Range("H46").Value = 1384
$ws.Range("I46").Value = 951.1111
$ws.Range("K46").Value = 951.1111
$ws.Range("M46").Value = -763.1111
$ws.Range("H55").Value = 841
$ws.Range("I55").Value = 205
$ws.Range("K55").Value = 205
$ws.Range("M55").Value = -32
$ws.Range("H68").Value = 2148.6
$ws.Range("I68").Value = 1872
$ws.Range("J68").Value = 2333
$ws.Range("K68").Value = 1872
$ws.Range("L68").Value = 2333
$ws.Range("M68").Value = -1123
$ws.Range("N68").Value = -3831
$ws.Range("H71").Value = 2148.6
$ws.Range("I71").Value = 1872
$ws.Range("J71").Value = 2333
$ws.Range("K71").Value = 9360
$ws.Range("L71").Value = 11665
$ws.Range("M71").Value = -5616
$ws.Range("N71").Value = -19153
$ws.Range("H88").Value = 32499.75
$ws.Range("I88").Value = 29999.666
$ws.Range("J88").Value = 40000
$ws.Range("K88").Value = 29999.666
$ws.Range("L88").Value = 40000
$ws.Range("M88").Value = -29571.666
$ws.Range("N88").Value = -40856
$ws.Range("H91").Value = 32499.75
$ws.Range("I91").Value = 29999.666
$ws.Range("J91").Value = 40000
$ws.Range("K91").Value = 29999.666
$ws.Range("L91").Value = 40000
$ws.Range("M91").Value = -28517.666
$ws.Range("N91").Value = -42964
$ws.Range("H126").Value = 2173
$ws.Range("I126").Value = 2201.8948
$ws.Range("K126").Value = 6605.6844
$ws.Range("M126").Value = -4135.6844
$ws.Range("H132").Value = 2051.1738
$ws.Range("I132").Value = 1998.3077
$ws.Range("J132").Value = 2119.9
$ws.Range("K132").Value = 5994.9231
$ws.Range("L132").Value = 6359.700000000001
$ws.Range("M132").Value = -3464.9231
$ws.Range("N132").Value = -11419.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 229.66667
$ws.Range("I107").Value = 229.66667
$ws.Range("K107").Value = 689.00001
$ws.Range("M107").Value = 1230.99999
$ws.Range("H122").Value = 2164.1333
$ws.Range("I122").Value = 2068.7144
$ws.Range("K122").Value = 6206.1432
$ws.Range("M122").Value = -3756.1432
$ws.Range("H126").Value = 3862.739
$ws.Range("I126").Value = 3313.0527
$ws.Range("J126").Value = 6473.75
$ws.Range("K126").Value = 9939.158100000001
$ws.Range("L126").Value = 19421.25
$ws.Range("M126").Value = -7469.158100000001
$ws.Range("N126").Value = -24361.25
$ws.Range("H130").Value = 37553.5
$ws.Range("J130").Value = 37553.5
$ws.Range("L130").Value = 37553.5
$ws.Range("N130").Value = -47593.5
$ws.Range("H132").Value = 3599.3794
$ws.Range("I132").Value = 3210.1155
$ws.Range("K132").Value = 9630.3465
$ws.Range("M132").Value = -7100.3465
